$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 69 through 108 in column C hold "TATA 1" -> change to "TATA 2"
for ($r = 69; $r -le 108; $r++) {
    $ws.Cells.Item($r, 3).Value = "TATA 2"
}

# Update the view state to match the saved selection/scroll position
$ws.Range("C67:C108").Select()
$excel.ActiveWindow.ScrollRow = 89
